$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel COM does not
# silently coerce them into Double values (which would drop formatting like
# trailing zeros, e.g. "1.000" -> 1).
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D22","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D34","D35","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.726.22"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "1.795.38"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "308.66"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "0.4406"
$ws.Range("E7").Value = "  +4.80%  "
$ws.Range("D8").Value = "0.3670"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "0.07331"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").Value = "0.8555"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").Value = "20.57"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value = "1.796.04"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "6.613"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "0.07060"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "91.24"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "5.255"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("D18").Value = "0.000008648"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "26.720.38"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("D22").Value = "5.150"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "10.79"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "1.971"
$ws.Range("D25").Value = "151.57"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "2.187"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").Value = "5.159"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").Value = "117.11"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "0.08781"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "0.7326"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").Value = "1.143"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("D33").Value = "2.903"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").Value = "4.422"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").Value = "1.002"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("D37").Value = "0.01946"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "0.05163"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("E39").Value = "  +3.80%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.957"
$ws.Range("E40").Value = "  -4.05%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.804"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").Value = "0.1676"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").Value = "0.5020"
$ws.Range("E43").Value = "  +6.32%  "
$ws.Range("D44").Value = "8.406"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").Value = "1.953"
$ws.Range("E45").Value = "  +4.08%  "
$ws.Range("D46").Value = "10.35"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").Value = "104.81"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "1.651"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").Value = "0.06281"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "0.9112"
$ws.Range("E51").Value = "  +1.40%  "
